# The commit adds one new weekly price record for "Coliflor" (Macroferia
# Regional de Talca) into the data table. The new record is inserted as
# row 267, pushing the existing rows 267-399 down to 268-400 (the sheet's
# used range grows from A1:R399 to A1:R400).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 267; rows 267-399 shift down to 268-400.
$ws.Rows("267:267").Insert()

# Populate the newly inserted row 267 with the new record's data.
$ws.Range("A267").Value = 5
$ws.Range("B267").Value = "Macroferia Regional de Talca"
$ws.Range("C267").Value = "Maule"
$ws.Range("D267").Value = 45016
$ws.Range("E267").Value = 7
$ws.Range("F267").Value = 100112008
$ws.Range("G267").Value = "Coliflor"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 3000
$ws.Range("K267").Value = 800
$ws.Range("L267").Value = 800
$ws.Range("M267").Value = 800
$ws.Range("N267").Value = '$/unidad'
$ws.Range("O267").Value = "Región del Maule"
$ws.Range("P267").Value = 800
$ws.Range("Q267").Value = 1
$ws.Range("R267").Value = "Hortaliza"
